# "Generate Report for Handback" - refresh the localization-status report
# after a successful handback: the status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the handback timestamps advance to the
# new generation time, and the stale "handback file is not the latest"
# error detail is cleared now that the file is up to date.
$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# 1. Overview sheet: E2/F2 status text change
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# 2. zh-cn sheet updates
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-24 04:47:14"
$zhcn.Range("P2").Value = ""

# 3. de-de sheet updates
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-24 04:47:21"
$dede.Range("P2").Value = ""

# Column widths widen for the status column (now holding longer text) and
# shrink for the now-empty Error Detail column, matching the committed layout
# (values chosen so the host's column-width grid lands closest to the target).
$overview.Columns.Item(5).ColumnWidth = 29.1666666666667
$overview.Columns.Item(6).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(3).ColumnWidth = 29.1666666666667
$zhcn.Columns.Item(16).ColumnWidth = 12.8333333333333
$dede.Columns.Item(3).ColumnWidth = 29.1666666666667
$dede.Columns.Item(16).ColumnWidth = 12.8333333333333

$wb.Save()
